$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C5").Value = 10872
$ws.Range("C6:C11").Value = 9774
$ws.Range("C12:C15").Value = 9541
$ws.Range("C16:C23").Value = 9371
$ws.Range("C24").Value = 9218
$ws.Range("C25:C30").Value = 9211
$ws.Range("C31:C44").Value = 8741
$ws.Range("C45:C61").Value = 8727
$ws.Range("C62:C65").Value = 8175
$ws.Range("C66:C67").Value = 7888
$ws.Range("C68:C99").Value = 7812
$ws.Range("C100:C252").Value = 7573
